# C1--C2-and-C3-PowerPoint.pptx edit
#
# 1) The table on slide 16 switches from table style
#    {8D10D1D9-8D14-466E-B808-38D4AD9B79F4} to
#    {0E89AE8D-F123-41B8-8F87-8485A1827FBC}.
# 2) The deck's "Integral" theme colours (theme2.xml, the theme actually
#    driving the slide master / design) are replaced with the stock
#    "Office Theme" colours (the palette that used to live, unused, in
#    theme1.xml). Table.Style can't be assigned directly - PowerPoint
#    requires Table.ApplyStyle(guid) - and individual theme colours are
#    changed through ThemeColorScheme.Colors(index).RGB.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{0E89AE8D-F123-41B8-8F87-8485A1827FBC}")

# --- 2. Swap the theme colour palette --------------------------------
# msoThemeColorDark1..msoThemeColorFollowedHyperlink (1..12) in order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# New values = the classic Office Theme palette.
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    $bgr = ($b * 0x10000) + ($g * 0x100) + $r
    $colorScheme.Colors($i).RGB = $bgr
}
